$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1100, 34.26159286499023, 3.047073125839233, 24.46402740478516, 46.90136337280273, 18476),
    @(1200, 34.30155944824219, 3.099591493606567, 24.6383056640625, 46.79884719848633, 18568),
    @(1300, 34.22674560546875, 3.08382248878479, 24.51870346069336, 45.56523895263672, 18438),
    @(1400, 34.08777618408203, 3.096240282058716, 24.39568328857422, 45.21668243408203, 18366),
    @(1500, 33.99085235595703, 3.170371055603027, 24.32392311096191, 48.96535110473633, 18392)
)

$startRow = 21

# Copy the Month column's cell format/type (text "09") down into the new
# rows first, so the new G cells stay text like G2:G20 (a plain .Value
# assignment of "09" would be auto-coerced to the number 9).
$ws.Range("G2").Copy()
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).PasteSpecial()
}

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
    $ws.Cells.Item($row, 5).Value = $data[$i][4]
    $ws.Cells.Item($row, 6).Value = $data[$i][5]
}
